# Applies the "first step towards #7" schema regen to the LinkML-generated
# Excel workbook:
#   * AssemblyJoinComponent + AssemblyJoin (two sheets) collapse into a
#     single AssemblyFragment sheet whose `location` column is split into
#     `left_location` / `right_location`.
#   * Every sheet with a repository_name/repository_id column pair gets the
#     two columns swapped (id first, then name), and the data validation
#     that lists the allowed repository names follows the renamed column.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# ---------------------------------------------------------------------
# 1. AssemblyJoinComponent -> AssemblyFragment: split `location` into
#    `left_location` / `right_location`.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("AssemblyJoinComponent")
$ws.Name = "AssemblyFragment"
$ws.Range("B1").Value = "left_location"
$ws.Range("C1").Value = "right_location"
$ws.Range("D1").Value = "reverse_complemented"

# ---------------------------------------------------------------------
# 2. AssemblyJoin sheet is obsolete now that AssemblyFragment carries both
#    edges directly - drop it.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("AssemblyJoin").Delete() | Out-Null

# ---------------------------------------------------------------------
# 3. RepositoryIdSource: repository_name/repository_id live in A1/B1.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("RepositoryIdSource")
$ws.Range("A2:A1048576").Validation.Delete() | Out-Null
$ws.Range("A1").Value = "repository_id"
$ws.Range("B1").Value = "repository_name"
$ws.Range("B2:B1048576").Validation.Add(3, 1, 1, '"addgene,genbank,benchling"') | Out-Null

# ---------------------------------------------------------------------
# 4. AddGeneIdSource: repository_name/repository_id live in C1/D1.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("AddGeneIdSource")
$ws.Range("C2:C1048576").Validation.Delete() | Out-Null
$ws.Range("C1").Value = "repository_id"
$ws.Range("D1").Value = "repository_name"
$ws.Range("D2:D1048576").Validation.Add(3, 1, 1, '"addgene,genbank,benchling"') | Out-Null

# ---------------------------------------------------------------------
# 5. BenchlingUrlSource: repository_name/repository_id live in A1/B1.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("BenchlingUrlSource")
$ws.Range("A2:A1048576").Validation.Delete() | Out-Null
$ws.Range("A1").Value = "repository_id"
$ws.Range("B1").Value = "repository_name"
$ws.Range("B2:B1048576").Validation.Add(3, 1, 1, '"addgene,genbank,benchling"') | Out-Null
